$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 7).Value = 1.72
$ws.Cells.Item(3, 8).Value = 3.2
$ws.Cells.Item(3, 9).Value = 5
$ws.Cells.Item(3, 12).Value = 1.5
$ws.Cells.Item(3, 13).Value = 2.27
$ws.Cells.Item(3, 14).Value = 2.4
$ws.Cells.Item(3, 15).Value = 1.44
$ws.Cells.Item(3, 16).Value = 1.55
$ws.Cells.Item(3, 17).Value = 2.15
$ws.Cells.Item(3, 21).Value = 6.6
$ws.Cells.Item(3, 22).Value = 9
$ws.Cells.Item(3, 23).Value = 13
$ws.Cells.Item(3, 24).Value = 18
$ws.Cells.Item(3, 26).Value = 5.3
$ws.Cells.Item(3, 27).Value = 6.6
$ws.Cells.Item(3, 28).Value = 24
$ws.Cells.Item(3, 30).Value = 9.75
$ws.Cells.Item(3, 31).Value = 26
$ws.Cells.Item(3, 32).Value = 18.5
$ws.Cells.Item(3, 33).Value = 100
$ws.Cells.Item(3, 34).Value = 70
$ws.Cells.Item(4, 9).Value = 3.5
$ws.Cells.Item(4, 12).Value = 1.33
$ws.Cells.Item(4, 13).Value = 2.77
$ws.Cells.Item(4, 14).Value = 1.98
$ws.Cells.Item(4, 15).Value = 1.65
$ws.Cells.Item(4, 16).Value = 1.44
$ws.Cells.Item(4, 17).Value = 2.42
$ws.Cells.Item(4, 18).Value = 1.83
$ws.Cells.Item(4, 19).Value = 1.78
$ws.Cells.Item(4, 20).Value = 6.7
$ws.Cells.Item(4, 21).Value = 8.75
$ws.Cells.Item(4, 24).Value = 16.5
$ws.Cells.Item(4, 25).Value = 30
$ws.Cells.Item(4, 26).Value = 9
$ws.Cells.Item(4, 28).Value = 16
$ws.Cells.Item(4, 29).Value = 80
$ws.Cells.Item(4, 30).Value = 9.5
$ws.Cells.Item(4, 31).Value = 18
$ws.Cells.Item(4, 34).Value = 35
$ws.Cells.Item(4, 35).Value = 45
$ws.Cells.Item(4, 36).Value = 700
$ws.Cells.Item(5, 7).Value = 2.6
$ws.Cells.Item(5, 8).Value = 2.88
$ws.Cells.Item(5, 9).Value = 2.88
$ws.Cells.Item(5, 10).Value = 1.11
$ws.Cells.Item(5, 11).Value = 6.5
$ws.Cells.Item(5, 12).Value = 1.44
$ws.Cells.Item(5, 13).Value = 2.63
$ws.Cells.Item(5, 20).Value = 7
$ws.Cells.Item(5, 24).Value = 26
$ws.Cells.Item(5, 27).Value = 5.5
$ws.Cells.Item(6, 7).Value = 1.36
$ws.Cells.Item(6, 8).Value = 5.25
$ws.Cells.Item(6, 10).Value = 1.03
$ws.Cells.Item(6, 11).Value = 17
$ws.Cells.Item(6, 12).Value = 1.2
$ws.Cells.Item(6, 13).Value = 4.33
$ws.Cells.Item(6, 14).Value = 1.65
$ws.Cells.Item(6, 15).Value = 2.2
$ws.Cells.Item(6, 16).Value = 1.29
$ws.Cells.Item(6, 17).Value = 3.5
$ws.Cells.Item(6, 18).Value = 1.91
$ws.Cells.Item(6, 19).Value = 1.8
$ws.Cells.Item(6, 20).Value = 7.5
$ws.Cells.Item(6, 21).Value = 7
$ws.Cells.Item(6, 22).Value = 8.5
$ws.Cells.Item(6, 24).Value = 11
$ws.Cells.Item(6, 25).Value = 26
$ws.Cells.Item(6, 26).Value = 15
$ws.Cells.Item(6, 27).Value = 11
$ws.Cells.Item(6, 28).Value = 21
$ws.Cells.Item(6, 29).Value = 67
$ws.Cells.Item(6, 30).Value = 17
$ws.Cells.Item(6, 35).Value = 41
$ws.Cells.Item(6, 36).Value = 351
$ws.Cells.Item(7, 7).Value = 3.6
$ws.Cells.Item(7, 8).Value = 3.4
$ws.Cells.Item(7, 9).Value = 1.95
$ws.Cells.Item(7, 12).Value = 1.29
$ws.Cells.Item(7, 13).Value = 3.5
$ws.Cells.Item(7, 18).Value = 1.8
$ws.Cells.Item(7, 19).Value = 1.91
$ws.Cells.Item(7, 21).Value = 19
$ws.Cells.Item(7, 24).Value = 29
$ws.Cells.Item(7, 27).Value = 6.5
$ws.Cells.Item(7, 30).Value = 7.5
$ws.Cells.Item(7, 31).Value = 9.5
$ws.Cells.Item(7, 32).Value = 9
$ws.Cells.Item(7, 33).Value = 17
$ws.Cells.Item(7, 34).Value = 17
$ws.Cells.Item(10, 8).Value = 5.5
$ws.Cells.Item(10, 9).Value = 12.5
$ws.Cells.Item(10, 11).Value = 8.5
$ws.Cells.Item(10, 12).Value = 1.21
$ws.Cells.Item(10, 13).Value = 3.9
$ws.Cells.Item(10, 14).Value = 1.65
$ws.Cells.Item(10, 15).Value = 2.1
$ws.Cells.Item(10, 16).Value = 1.33
$ws.Cells.Item(10, 17).Value = 3.05
$ws.Cells.Item(10, 18).Value = 2.37
$ws.Cells.Item(10, 19).Value = 1.52
$ws.Cells.Item(10, 20).Value = 6.5
$ws.Cells.Item(10, 24).Value = 11.75
$ws.Cells.Item(10, 26).Value = 8.5
$ws.Cells.Item(10, 27).Value = 11.5
$ws.Cells.Item(10, 29).Value = 175
$ws.Cells.Item(10, 31).Value = 100
$ws.Cells.Item(10, 32).Value = 40
$ws.Cells.Item(10, 33).Value = 500
$ws.Cells.Item(10, 34).Value = 200
$ws.Cells.Item(10, 35).Value = 150
$ws.Cells.Item(13, 7).Value = 4.05
$ws.Cells.Item(13, 8).Value = 2.87
$ws.Cells.Item(13, 9).Value = 2.05
$ws.Cells.Item(13, 11).Value = 6.5
$ws.Cells.Item(13, 12).Value = 1.33
$ws.Cells.Item(13, 13).Value = 3.05
$ws.Cells.Item(13, 14).Value = 1.98
$ws.Cells.Item(13, 15).Value = 1.75
$ws.Cells.Item(13, 16).Value = 1.38
$ws.Cells.Item(13, 17).Value = 2.82
$ws.Cells.Item(13, 20).Value = 11.25
$ws.Cells.Item(13, 21).Value = 24
$ws.Cells.Item(13, 26).Value = 6.5
$ws.Cells.Item(13, 27).Value = 5.6
$ws.Cells.Item(13, 31).Value = 10.25
$ws.Cells.Item(13, 32).Value = 8.25
$ws.Cells.Item(14, 7).Value = 2.2
$ws.Cells.Item(14, 8).Value = 3.3
$ws.Cells.Item(14, 9).Value = 2.85
$ws.Cells.Item(14, 22).Value = 7.5
$ws.Cells.Item(14, 23).Value = 18
$ws.Cells.Item(14, 24).Value = 13.5
$ws.Cells.Item(14, 25).Value = 18
$ws.Cells.Item(14, 26).Value = 12
$ws.Cells.Item(14, 27).Value = 5.9
$ws.Cells.Item(14, 28).Value = 10
$ws.Cells.Item(14, 29).Value = 35
$ws.Cells.Item(14, 30).Value = 9.25
$ws.Cells.Item(14, 32).Value = 8.75
$ws.Cells.Item(14, 34).Value = 18
$ws.Cells.Item(14, 35).Value = 21
$ws.Cells.Item(14, 36).Value = 200
$ws.Cells.Item(15, 7).Value = 1.19
$ws.Cells.Item(15, 8).Value = 5
$ws.Cells.Item(15, 9).Value = 14.5
$ws.Cells.Item(15, 14).Value = 1.55
$ws.Cells.Item(15, 15).Value = 2.15
$ws.Cells.Item(15, 21).Value = 4.9
$ws.Cells.Item(15, 22).Value = 7.8
$ws.Cells.Item(15, 23).Value = 5.7
$ws.Cells.Item(15, 24).Value = 9
$ws.Cells.Item(15, 25).Value = 26
$ws.Cells.Item(15, 26).Value = 11.75
$ws.Cells.Item(15, 27).Value = 9.5
$ws.Cells.Item(15, 28).Value = 23
$ws.Cells.Item(15, 29).Value = 110
$ws.Cells.Item(15, 31).Value = 110
$ws.Cells.Item(15, 32).Value = 40
$ws.Cells.Item(17, 23).Value = 18.5
$ws.Cells.Item(17, 24).Value = 17
$ws.Cells.Item(17, 29).Value = 75
$ws.Cells.Item(17, 30).Value = 10
$ws.Cells.Item(17, 31).Value = 20
$ws.Cells.Item(18, 20).Value = 6.5
$ws.Cells.Item(18, 21).Value = 12
$ws.Cells.Item(18, 22).Value = 10.5
$ws.Cells.Item(18, 24).Value = 28
$ws.Cells.Item(18, 30).Value = 6.8
$ws.Cells.Item(18, 32).Value = 10.5
$ws.Cells.Item(18, 33).Value = 32
$ws.Cells.Item(18, 34).Value = 28
$ws.Cells.Item(18, 35).Value = 45
$ws.Cells.Item(19, 12).Value = 1.25
$ws.Cells.Item(19, 13).Value = 3.75
$ws.Cells.Item(19, 14).Value = 1.88
$ws.Cells.Item(19, 15).Value = 1.93
$ws.Cells.Item(21, 20).Value = 19.5
$ws.Cells.Item(24, 10).Value = 1.07
$ws.Cells.Item(24, 12).Value = 1.36
$ws.Cells.Item(24, 18).Value = 1.91
$ws.Cells.Item(24, 19).Value = 1.91
$ws.Cells.Item(25, 7).Value = 1.45
$ws.Cells.Item(25, 8).Value = 5
$ws.Cells.Item(25, 9).Value = 5.75
$ws.Cells.Item(25, 21).Value = 10
$ws.Cells.Item(25, 22).Value = 9
$ws.Cells.Item(25, 23).Value = 12
$ws.Cells.Item(25, 31).Value = 41
$ws.Cells.Item(25, 32).Value = 19
$ws.Cells.Item(25, 33).Value = 67
$ws.Cells.Item(25, 34).Value = 41
$ws.Cells.Item(25, 35).Value = 34
$ws.Cells.Item(26, 3).Value = "11:00"
$ws.Cells.Item(28, 30).Value = 6.3
$ws.Cells.Item(31, 10).Value = 1.07
$ws.Cells.Item(31, 11).Value = 9
$ws.Cells.Item(31, 12).Value = 1.36
$ws.Cells.Item(31, 13).Value = 3
$ws.Cells.Item(31, 18).Value = 1.83
$ws.Cells.Item(31, 19).Value = 1.83
$ws.Cells.Item(32, 7).Value = 3.5
$ws.Cells.Item(32, 8).Value = 3.25
$ws.Cells.Item(32, 10).Value = 1.08
$ws.Cells.Item(32, 11).Value = 7.5
$ws.Cells.Item(32, 12).Value = 1.44
$ws.Cells.Item(32, 13).Value = 2.63
$ws.Cells.Item(32, 14).Value = 2.35
$ws.Cells.Item(32, 15).Value = 1.57
$ws.Cells.Item(32, 21).Value = 17
$ws.Cells.Item(32, 22).Value = 13
$ws.Cells.Item(32, 24).Value = 34
$ws.Cells.Item(32, 26).Value = 7.5
$ws.Cells.Item(32, 27).Value = 6
$ws.Cells.Item(32, 31).Value = 9
$ws.Cells.Item(32, 33).Value = 19
$ws.Cells.Item(32, 34).Value = 21
$ws.Cells.Item(32, 36).Value = 451
$ws.Cells.Item(34, 18).Value = 1.53
$ws.Cells.Item(34, 19).Value = 2.38
$ws.Cells.Item(35, 14).Value = 1.82
$ws.Cells.Item(35, 15).Value = 1.87
$ws.Cells.Item(36, 16).Value = 1.27
$ws.Cells.Item(37, 18).Value = 1.67
$ws.Cells.Item(38, 10).Value = 1.04
$ws.Cells.Item(38, 12).Value = 1.25
$ws.Cells.Item(38, 14).Value = 1.83
$ws.Cells.Item(38, 15).Value = 1.98
$ws.Cells.Item(39, 10).Value = 1.04
$ws.Cells.Item(39, 12).Value = 1.22
$ws.Cells.Item(39, 18).Value = 1.7
$ws.Cells.Item(41, 16).Value = 1.27
$ws.Cells.Item(44, 18).Value = 1.44
$ws.Cells.Item(44, 19).Value = 2.63
$ws.Cells.Item(45, 12).Value = 1.17
$ws.Cells.Item(45, 13).Value = 5
$ws.Cells.Item(45, 14).Value = 1.57
$ws.Cells.Item(45, 15).Value = 2.35
